$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 30; existing rows 30-41 shift down to 31-42.
$ws.Rows.Item(30).Insert()

# Populate the newly inserted row 30 with the new Membrillo price record.
$ws.Cells.Item(30, 1).Value = 5
$ws.Cells.Item(30, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(30, 3).Value = "Maule"
$ws.Cells.Item(30, 4).Value = 44680
$ws.Cells.Item(30, 5).Value = 7
$ws.Cells.Item(30, 6).Value = "Fruta"
$ws.Cells.Item(30, 7).Value = 100104
$ws.Cells.Item(30, 8).Value = "Frutos de pepita"
$ws.Cells.Item(30, 9).Value = 100104003
$ws.Cells.Item(30, 10).Value = "Membrillo"
$ws.Cells.Item(30, 11).Value = "Champion"
$ws.Cells.Item(30, 12).Value = "Primera"
$ws.Cells.Item(30, 13).Value = 230
$ws.Cells.Item(30, 14).Value = 10000
$ws.Cells.Item(30, 15).Value = 10000
$ws.Cells.Item(30, 16).Value = 10000
$ws.Cells.Item(30, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(30, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(30, 19).Value = 556
$ws.Cells.Item(30, 20).Value = 18

# Keep the date column formatted like the rest of column D.
$ws.Cells.Item(30, 4).NumberFormat = $ws.Cells.Item(31, 4).NumberFormat
